$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a NIB-like value (leading-zero numeric string) as true text,
# preserving leading zeros, then restore the default "Normal" cell style so
# no stray direct formatting is left on the cell.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Rows 22-32: newly added data rows (default style - no direct formatting)
$rows = @(
    @{ r = 22; a = "HUSNAWATI,S.PD.-";              b = "00632"; c = 20000 },
    @{ r = 23; a = "HJ. SITTI HARIANA, S.ST,KER";    b = "00095"; c = 9207  },
    @{ r = 24; a = "SUJALIL, ST";                    b = "00145"; c = 4237  },
    @{ r = 25; a = "LA FAIDI, SP";                   b = "00144"; c = 6480  },
    @{ r = 26; a = "IDHAM KASMIN, SE";               b = "00257"; c = 9665  },
    @{ r = 27; a = "IKHVAN KASMIN, SH";              b = "00259"; c = 9665  },
    @{ r = 28; a = "TUTI ALAMIA, S. PD SD";           b = "00260"; c = 999   },
    @{ r = 29; a = "ARDIN, SE";                      b = "00379"; c = 703   },
    @{ r = 30; a = "TARFAN, S. PD";                  b = "00316"; c = 11790 },
    @{ r = 31; a = "TARFAN, S. PD";                  b = "00398"; c = 915   },
    @{ r = 32; a = "TARFAN, S. PD";                  b = "00317"; c = 12310 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    Set-TextValue $ws.Range("B$r") $row.b
    $ws.Range("C$r").Value = $row.c
}

# Row 33: reuses existing values, styled like the other body rows (style of A4/B4/C4)
$ws.Range("A33").Value = "HAYUN ANWAR"
Set-TextValue $ws.Range("B33") "00218"
$ws.Range("C33").Value = 20000
$ws.Range("A33:C33").Style = $ws.Range("A4:C4").Style

# Row 21 (A21/C21) loses its stray "applyFill" direct style, matching the
# plain body-row style used elsewhere (same visible style as A2/C2).
$ws.Range("A21").Style = $ws.Range("A2").Style
$ws.Range("C21").Style = $ws.Range("C2").Style

# Reflect the final selection left by the author after entering the new row.
$ws.Activate()
$ws.Range("A33:C33").Select()
